$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-converted to numbers by Excel (values that look like plain numbers).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'

$ws.Range('D2').Value = '66.191.96'
$ws.Range('E2').Value = '  -1.05%  '
$ws.Range('D3').Value = '3.275.10'
$ws.Range('E3').Value = '  -1.52%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '585.04'
$ws.Range('E5').Value = '  +1.66%  '
$ws.Range('D6').Value = '179.06'
$ws.Range('E6').Value = '  -2.11%  '
$ws.Range('E7').Value = '  +3.94%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '0.125'
$ws.Range('D10').Value = '6.73'
$ws.Range('E10').Value = '  +1.43%  '
$ws.Range('E11').Value = '  -0.93%  '
$ws.Range('D12').Value = '3.846.83'
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('E13').Value = '  -3.56%  '
$ws.Range('D14').Value = '66.208.72'
$ws.Range('E14').Value = '  -1.37%  '
$ws.Range('D15').Value = '26.32'
$ws.Range('E15').Value = '  -3.13%  '
$ws.Range('E16').Value = '  -1.84%  '
$ws.Range('D17').Value = '3.282.66'
$ws.Range('E17').Value = '  -1.25%  '
$ws.Range('D18').Value = '432.98'
$ws.Range('E18').Value = '  -2.10%  '
$ws.Range('D19').Value = '5.50'
$ws.Range('E19').Value = '  -2.80%  '
$ws.Range('D20').Value = '13.15'
$ws.Range('E20').Value = '  -2.93%  '
$ws.Range('D21').Value = '7.38'
$ws.Range('E21').Value = '  -4.15%  '
$ws.Range('D22').Value = '71.63'
$ws.Range('E22').Value = '  -3.09%  '
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').Value = '3.423.55'
$ws.Range('E24').Value = '  -1.55%  '
$ws.Range('E25').Value = '  -1.34%  '
$ws.Range('E26').Value = '  +0.65%  '
$ws.Range('E27').Value = '  -5.06%  '
$ws.Range('D28').Value = '8.82'
$ws.Range('E28').Value = '  -1.67%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  -0.90%  '
$ws.Range('D31').Value = '22.23'
$ws.Range('E31').Value = '  -2.99%  '
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('D33').Value = '5.15'
$ws.Range('E33').Value = '  -2.74%  '
$ws.Range('D34').Value = '6.59'
$ws.Range('E34').Value = '  -2.85%  '
$ws.Range('E35').Value = '  -2.79%  '
$ws.Range('D36').Value = '157.85'
$ws.Range('E36').Value = '  -2.69%  '
$ws.Range('D37').Value = '1.42'
$ws.Range('E37').Value = '  -5.06%  '
$ws.Range('D38').Value = '26.49'
$ws.Range('E38').Value = '  -3.92%  '
$ws.Range('D39').Value = '1.78'
$ws.Range('E39').Value = '  -2.96%  '
$ws.Range('D40').Value = '2.773.99'
$ws.Range('E40').Value = '  -1.86%  '
$ws.Range('E41').Value = '  -2.03%  '
$ws.Range('E42').Value = '  -2.77%  '
$ws.Range('E43').Value = '  -0.41%  '
$ws.Range('E44').Value = '  -2.69%  '
$ws.Range('E45').Value = '  -2.16%  '
$ws.Range('D46').Value = '320.87'
$ws.Range('E46').Value = '  -0.24%  '
$ws.Range('E47').Value = '  -0.64%  '
$ws.Range('D48').Value = '23.17'
$ws.Range('E48').Value = '  -4.60%  '
$ws.Range('D49').Value = '0.0266'
$ws.Range('E49').Value = '  -2.30%  '
$ws.Range('E50').Value = '  +2.83%  '
$ws.Range('E51').Value = '  -0.01%  '
